$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, pushing existing rows 29..39 down to 30..40
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new weekly record
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C29").Value = 'Coquimbo'
$ws.Range("D29").Value = 44876
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112026
$ws.Range("G29").Value = 'Haba'
$ws.Range("H29").Value = 'Sin especificar'
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 460
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 6500
$ws.Range("N29").Value = '$/saco 25 kilos'
$ws.Range("O29").Value = 'Provincia del Elquí'
$ws.Range("P29").Value = 260
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = 'Hortaliza'
